$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "ID",
    "Seal of State",
    "Seal of State Issuer",
    "Second Seal",
    "Second Seal Issuer",
    "Bueraucratic Stamp",
    "Place of Studio's Photographer's Name",
    "Photographer",
    "Location of Photographer",
    "Date of Document",
    "Date on Photograph",
    "Handwritten on front",
    "Numbered",
    "Perforated",
    "Printed information on Front",
    "Writing on Front",
    "Date of Photograph",
    "Color of Ink",
    "Other notes"
)

$col = 2
foreach ($header in $headers) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $header
    $col = $col + 1
}

$range = $ws.Range("B1:T1")
$range.Font.Bold = $true
$range.HorizontalAlignment = -4108  # xlCenter
$range.VerticalAlignment = -4160    # xlTop
$range.Borders.LineStyle = 1        # xlContinuous
$range.Borders.Weight = 2           # xlThin
$range.Borders.ColorIndex = 1       # automatic
